$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate formatting (styles) of the last existing data row onto the two new rows
$ws.Range("A163:V163").Copy()
$ws.Range("A164:V164").PasteSpecial(-4122)
$ws.Range("A165:V165").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 164
$ws.Range("A164").Value = 163
$ws.Range("B164").Value = "romania"
$ws.Range("C164").Value = "liga-1"
$ws.Range("D164").Value = "2023-2024"
$ws.Range("E164").Value = 45280.6875
$ws.Range("F164").Value = "UTA Arad"
$ws.Range("G164").Value = 2
$ws.Range("H164").Value = "Univ. Craiova"
$ws.Range("I164").Value = 2
$ws.Range("J164").Value = 3.95
$ws.Range("K164").Value = "16/12/2023 18:13"
$ws.Range("L164").Value = 5.24
$ws.Range("M164").Value = "20/12/2023 16:24"
$ws.Range("N164").Value = 3.44
$ws.Range("O164").Value = "16/12/2023 18:13"
$ws.Range("P164").Value = 3.44
$ws.Range("Q164").Value = "20/12/2023 16:24"
$ws.Range("R164").Value = 1.87
$ws.Range("S164").Value = "16/12/2023 18:13"
$ws.Range("T164").Value = 1.76
$ws.Range("U164").Value = "20/12/2023 16:24"
$ws.Range("V164").Value = "https://www.betexplorer.com/football/romania/liga-1/fc-uta-arad-univ-craiova/IHt8EA9a/"

# Row 165
$ws.Range("A165").Value = 164
$ws.Range("B165").Value = "romania"
$ws.Range("C165").Value = "liga-1"
$ws.Range("D165").Value = "2023-2024"
$ws.Range("E165").Value = 45280.8125
$ws.Range("F165").Value = "Farul Constanta"
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = "FC Rapid Bucuresti"
$ws.Range("I165").Value = 0
$ws.Range("J165").Value = 2.3
$ws.Range("K165").Value = "16/12/2023 13:12"
$ws.Range("L165").Value = 2.24
$ws.Range("M165").Value = "20/12/2023 19:21"
$ws.Range("N165").Value = 3.28
$ws.Range("O165").Value = "16/12/2023 13:12"
$ws.Range("P165").Value = 3.39
$ws.Range("Q165").Value = "20/12/2023 19:21"
$ws.Range("R165").Value = 2.98
$ws.Range("S165").Value = "16/12/2023 13:12"
$ws.Range("T165").Value = 3.31
$ws.Range("U165").Value = "20/12/2023 19:21"
$ws.Range("V165").Value = "https://www.betexplorer.com/football/romania/liga-1/farul-constanta-rapid-bucuresti/CvtCDUO5/"
